$d = $word.ActiveDocument

# Turn on revision tracking so the new sentence is inserted as its own
# run (rather than being silently merged back into the run it was typed
# into) - this mirrors how Word keeps freshly-typed text in a distinct
# run from its neighbours. We accept the revision immediately afterwards
# so the final document contains plain (non-tracked) runs.
$d.TrackRevisions = $true

$rng = $d.Content
$rng.Find.Execute("localhost:8080.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Collapse(0)
$rng.InsertAfter(" This is most easily achieved by building and running the docker image from the supplied DOCKERFILE.")

$d.TrackRevisions = $false
$d.AcceptAllRevisions()
